# Update countries & provincias Spain
# Applies the 23-Oct-2020 refresh: updated case counts for several
# countries, a handful of rows whose country label shifted by one
# position (because the underlying source re-sorted when new rows were
# merged in), and a refreshed "last updated" timestamp in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 12:05"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("B4").Value = 8664365
$ws.Range("C4").Value = 2714
$ws.Range("D4").Value = 5656150
$ws.Range("E4").Value = 2779792
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 228423

# --- Row 5: India ------------------------------------------------------------
$ws.Range("B5").Value = 7763067
$ws.Range("C5").Value = 3427
$ws.Range("E5").Value = 697205
$ws.Range("G5").Value = 29
$ws.Range("H5").Value = 117365

# --- Row 20: Alemania --------------------------------------------------------
$ws.Range("B20").Value = 404189
$ws.Range("C20").Value = 315
$ws.Range("E20").Value = 88045

# --- Row 22: Indonesia --------------------------------------------------------
$ws.Range("B22").Value = 381910
$ws.Range("C22").Value = 4369
$ws.Range("D22").Value = 305100
$ws.Range("E22").Value = 63733
$ws.Range("G22").Value = 118
$ws.Range("H22").Value = 13077

# --- Row 51: Suiza -------------------------------------------------------------
$ws.Range("E51").Value = 39263
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 2056

# --- Rows 80-83: Eslovaquia moves ahead of Dinamarca / Serbia / Bosnia y Herzegovina
# (new country order: Afganistan, Eslovaquia, Dinamarca, Serbia, Bosnia y Herzegovina, Bulgaria)
$ws.Range("A80").Value = "Eslovaquia"
$ws.Range("B80").Value = 37911
$ws.Range("C80").Value = 2581
$ws.Range("D80").Value = 8859
$ws.Range("E80").Value = 28918
$ws.Range("G80").Value = 19
$ws.Range("H80").Value = 134

$ws.Range("A81").Value = "Dinamarca"
$ws.Range("B81").Value = 37763
$ws.Range("D81").Value = 30877
$ws.Range("E81").Value = 6192
$ws.Range("H81").Value = 694

$ws.Range("A82").Value = "Serbia"
$ws.Range("B82").Value = 37536
$ws.Range("D82").Value = 31536
$ws.Range("E82").Value = 5217
$ws.Range("H82").Value = 783

$ws.Range("A83").Value = "Bosnia y Herzegovina"
$ws.Range("B83").Value = 37314
$ws.Range("D83").Value = 25989
$ws.Range("E83").Value = 10274
$ws.Range("H83").Value = 1051

# --- Row 85: El Salvador -------------------------------------------------------
$ws.Range("B85").Value = 32421
$ws.Range("C85").Value = 159
$ws.Range("D85").Value = 28127
$ws.Range("E85").Value = 3354

# --- Row 102: Finlandia ---------------------------------------------------------
$ws.Range("B102").Value = 14474
$ws.Range("C102").Value = 219
$ws.Range("E102").Value = 4319

# --- Row 123: Sri Lanka -----------------------------------------------------------
$ws.Range("D123").Value = 3644
$ws.Range("E123").Value = 2629

# --- Row 131: Hong Kong ------------------------------------------------------------
$ws.Range("B131").Value = 5285
$ws.Range("C131").Value = 4
$ws.Range("D131").Value = 5029
$ws.Range("E131").Value = 151

# --- Rows 144-145: Letonia moves ahead of Mayotte ------------------------------------
$ws.Range("A144").Value = "Letonia"
$ws.Range("B144").Value = 4208
$ws.Range("C144").Value = 250
$ws.Range("D144").Value = 1357
$ws.Range("E144").Value = 2801
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 50

$ws.Range("A145").Value = "Mayotte"
$ws.Range("B145").Value = 4203
$ws.Range("D145").Value = 2964
$ws.Range("E145").Value = 1195
$ws.Range("H145").Value = 44

# --- Rows 216-217: Montserrat moves ahead of Islas Malvinas --------------------------
$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

$ws.Range("A217").Value = "Islas Malvinas"
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0
